# "Überarbeitung / Anpassung Evaluation View"
# Rework of the "Gruppendaten" template header row: drop the competition
# columns (Gruppe/StartNR/Platz/Gesamtpunkte), move Geburtsdatum up next to
# the other personal-data fields, and append the new camp-logistics columns
# (LagerNr/Status/Essgewohnheiten/Unverträglichkeiten).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2) ------------------------------------------------
# New column order: OU | Feuerwehr | LagerNr | Geschlecht | Vorname |
#                    Nachname | Geburtsdatum | Alter | Status |
#                    Essgewohnheiten | Unverträglichkeiten
$ws.Range("A2").Value = "OU"
$ws.Range("B2").Value = "Feuerwehr"
$ws.Range("C2").Value = "LagerNr"
$ws.Range("D2").Value = "Geschlecht"
$ws.Range("E2").Value = "Vorname"
$ws.Range("F2").Value = "Nachname"
$ws.Range("G2").Value = "Geburtsdatum"
$ws.Range("H2").Value = "Alter"
$ws.Range("I2").Value = "Status"
$ws.Range("J2").Value = "Essgewohnheiten"
$ws.Range("K2").Value = "Unverträglichkeiten"

# A1 keeps its text ("Gruppendaten") - left untouched on purpose.

# --- Row heights ---------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 46.5
$ws.Rows.Item(2).RowHeight = 18.75

# --- Column widths (character units) -------------------------------------
# Re-sized to fit the new headers/content.
$ws.Columns.Item(1).ColumnWidth = 22.736979166666668
$ws.Columns.Item(2).ColumnWidth = 17.877604166666668
$ws.Columns.Item(3).ColumnWidth = 14.451822916666666
$ws.Columns.Item(4).ColumnWidth = 13.022135416666666
$ws.Columns.Item(5).ColumnWidth = 17.877604166666668
$ws.Columns.Item(6).ColumnWidth = 21.307291666666668
$ws.Columns.Item(7).ColumnWidth = 17.451822916666668
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 20.736979166666668
$ws.Columns.Item(11).ColumnWidth = 23.451822916666668

# --- Selection -------------------------------------------------------------
$ws.Range("K3").Select() | Out-Null
